# Rename header row labels on the "Data" sheet from Title Case display
# names to lowercase machine-style field names, and apply numeric
# formatting to the Year / Price columns.

$wb = $excel.ActiveWorkbook

# The "Quantity" sheet's A2 carried a stray bold-ish style left over from
# past editing; the refreshed workbook drops it back to the default style.
$wsQty = $wb.Worksheets.Item("Quantity")
$wsQty.Range("A2").Style = "Normal"

$ws = $wb.Worksheets.Item("Data")

$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "brand"
$ws.Range("C1").Value = "category"
$ws.Range("D1").Value = "colour"
$ws.Range("E1").Value = "frameSize"
$ws.Range("F1").Value = "frameMaterial"
$ws.Range("G1").Value = "year"
$ws.Range("H1").Value = "price"
$ws.Range("I1").Value = "image"

$ws.Range("G2:G41").NumberFormat = "0"
$ws.Range("H2:H41").NumberFormat = "0.00"

$ws.Range("E32").Select()
